$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.716.31"
$ws.Range("E2").Value = "  +4.56%  "
$ws.Range("D3").Value = "3.332.48"
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.51"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "3.908.47"
$ws.Range("E12").Value = "  +4.57%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.70"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").Value = "62.747.11"
$ws.Range("E16").Value = "  +4.68%  "
$ws.Range("D17").Value = "3.336.59"
$ws.Range("E17").Value = "  +5.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.50"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.91%  "
$ws.Range("E19").Value = "  +4.67%  "
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.38"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.61"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.177"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "0.0₃0956"
$ws.Range("E27").Value = "  +5.93%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.46"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.29%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.97"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.96"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.56"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.28"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.66"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.71%  "
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("E36").Value = "  +9.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +11.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.36"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.04%  "
$ws.Range("D39").Value = "2.851.49"
$ws.Range("E39").Value = "  +2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0734"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("E41").Value = "  +8.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.31"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.71"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.747"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.03"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("D46").Value = "3.380.30"
$ws.Range("E46").Value = "  +4.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.90"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.32%  "
$ws.Range("E48").Value = "  +3.23%  "
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.804"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "283.12"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.61%  "
